$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.0292345
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.4428765120700495
$ws.Range("J2").Value = 0.346386487911515
$ws.Range("M2").Value = 1.532593
$ws.Range("N2").Value = 3.065186
$ws.Range("O2").Value = 0.0795983245703594
$ws.Range("P2").Value = 0.05838920196386116
$ws.Range("Q2").Value = 0.0448045900585
$ws.Range("R2").Value = 0.179218360234
$ws.Range("S2").Value = 0.03525222835234049
$ws.Range("T2").Value = 0.020225230600218

# Row 3
$ws.Range("G3").Value = 0.0292345
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.4428765120700495
$ws.Range("J3").Value = 0.346386487911515
$ws.Range("O3").Value = 0.6326044366842063
$ws.Range("P3").Value = 0.6960687002426557
$ws.Range("Q3").Value = 0.3560826513348333
$ws.Range("R3").Value = 2.136495908009
$ws.Range("S3").Value = 0.2801656464387398
$ws.Range("T3").Value = 0.2411087924221866

# Row 4
$ws.Range("G4").Value = 0.0292345
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.4428765120700495
$ws.Range("J4").Value = 0.346386487911515
$ws.Range("M4").Value = 0.8528209999999999
$ws.Range("N4").Value = 2.558463
$ws.Range("O4").Value = 0.04429298760885536
$ws.Range("P4").Value = 0.04873655720209673
$ws.Range("Q4").Value = 0.0249317955245
$ws.Range("R4").Value = 0.149590773147
$ws.Range("S4").Value = 0.01961632386137178
$ws.Range("T4").Value = 0.01688168488213294

# Row 5
$ws.Range("G5").Value = 0.0292345
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.4428765120700495
$ws.Range("J5").Value = 0.346386487911515
$ws.Range("M5").Value = 3.7338975
$ws.Range("N5").Value = 7.467795
$ws.Range("O5").Value = 0.1939275366111247
$ws.Range("P5").Value = 0.142255181408147
$ws.Range("Q5").Value = 0.10915862646375
$ws.Range("R5").Value = 0.436634505855
$ws.Range("S5").Value = 0.08588595100867175
$ws.Range("T5").Value = 0.04927527267518349

# Row 6
$ws.Range("G6").Value = 0.0292345
$ws.Range("H6").Value = 0.058469
$ws.Range("I6").Value = 0.4428765120700495
$ws.Range("J6").Value = 0.346386487911515
$ws.Range("M6").Value = 0.2147316666666667
$ws.Range("N6").Value = 0.644195
$ws.Range("O6").Value = 0.01115252444639089
$ws.Range("P6").Value = 0.01227137014168456
$ws.Range("Q6").Value = 0.006277572909166666
$ws.Range("R6").Value = 0.037665437455
$ws.Range("S6").Value = 0.004939191127593557
$ws.Range("T6").Value = 0.004250636805240345

# Row 7
$ws.Range("G7").Value = 0.0292345
$ws.Range("H7").Value = 0.058469
$ws.Range("I7").Value = 0.4428765120700495
$ws.Range("J7").Value = 0.346386487911515
$ws.Range("M7").Value = 0.7398226666666666
$ws.Range("N7").Value = 2.219468
$ws.Range("O7").Value = 0.03842419007906348
$ws.Range("P7").Value = 0.04227898904155473
$ws.Range("Q7").Value = 0.02162834574866667
$ws.Range("R7").Value = 0.129770074492
$ws.Range("S7").Value = 0.01701717128133223
$ws.Range("T7").Value = 0.01464487052655357

# Row 8
$ws.Range("I8").Value = 0.5571234879299505
$ws.Range("J8").Value = 0.6536135120884849
$ws.Range("M8").Value = 1.532593
$ws.Range("N8").Value = 3.065186
$ws.Range("O8").Value = 0.0795983245703594
$ws.Range("P8").Value = 0.05838920196386116
$ws.Range("Q8").Value = 0.05636264016799999
$ws.Range("R8").Value = 0.338175841008
$ws.Range("S8").Value = 0.0443460962180189
$ws.Range("T8").Value = 0.03816397136364316

# Row 9
$ws.Range("I9").Value = 0.5571234879299505
$ws.Range("J9").Value = 0.6536135120884849
$ws.Range("O9").Value = 0.6326044366842063
$ws.Range("P9").Value = 0.6960687002426557
$ws.Range("S9").Value = 0.3524387902454665
$ws.Range("T9").Value = 0.454959907820469

# Row 10
$ws.Range("I10").Value = 0.5571234879299505
$ws.Range("J10").Value = 0.6536135120884849
$ws.Range("M10").Value = 0.8528209999999999
$ws.Range("N10").Value = 2.558463
$ws.Range("O10").Value = 0.04429298760885536
$ws.Range("P10").Value = 0.04873655720209673
$ws.Range("Q10").Value = 0.03136334509599999
$ws.Range("R10").Value = 0.282270105864
$ws.Range("S10").Value = 0.02467666374748358
$ws.Range("T10").Value = 0.03185487231996379

# Row 11
$ws.Range("I11").Value = 0.5571234879299505
$ws.Range("J11").Value = 0.6536135120884849
$ws.Range("M11").Value = 3.7338975
$ws.Range("N11").Value = 7.467795
$ws.Range("O11").Value = 0.1939275366111247
$ws.Range("P11").Value = 0.142255181408147
$ws.Range("Q11").Value = 0.13731781446
$ws.Range("R11").Value = 0.8239068867599999
$ws.Range("S11").Value = 0.108041585602453
$ws.Range("T11").Value = 0.09297990873296352

# Row 12
$ws.Range("I12").Value = 0.5571234879299505
$ws.Range("J12").Value = 0.6536135120884849
$ws.Range("M12").Value = 0.2147316666666667
$ws.Range("N12").Value = 0.644195
$ws.Range("O12").Value = 0.01115252444639089
$ws.Range("P12").Value = 0.01227137014168456
$ws.Range("Q12").Value = 0.007896971773333332
$ws.Range("R12").Value = 0.07107274595999999
$ws.Range("S12").Value = 0.006213333318797334
$ws.Range("T12").Value = 0.008020733336444215

# Row 13
$ws.Range("I13").Value = 0.5571234879299505
$ws.Range("J13").Value = 0.6536135120884849
$ws.Range("M13").Value = 0.7398226666666666
$ws.Range("N13").Value = 2.219468
$ws.Range("O13").Value = 0.03842419007906348
$ws.Range("P13").Value = 0.04227898904155473
$ws.Range("Q13").Value = 0.02720771838933333
$ws.Range("R13").Value = 0.244869465504
$ws.Range("S13").Value = 0.02140701879773125
$ws.Range("T13").Value = 0.02763411851500116
